# Adds Uncertainty to Model:
#  - new defined names UNC_Fluid_C (B9) and UNC_Fluid_Temp (B10)
#  - B3 fluid_mass value updated
#  - B7 Q_toBoil formula updated to include the new uncertainty factors, with a "J" units label in C7
#  - new rows 9/10 holding the uncertainty inputs + labels/units
#  - column A widened to fit the new longer labels
#  - selection moved to B15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update fluid_mass value
$ws.Range("B3").Value = 1.4137154999999999

# New uncertainty input rows (fill labels/units first so new shared-string
# entries land in the same order as the authored workbook)
$ws.Range("A9").Value = "UNC_Fluid_C"
$ws.Range("C9").Value = "C value multiplier for fluids other than water"
$ws.Range("B9").Value = 1.08

# Update the Q_toBoil formula to incorporate the new uncertainty terms
$ws.Range("B7").Formula = "=fluid_mass*c_water*UNC_Fluid_C*(set_temp-ambient_temp+UNC_Fluid_Temp)"
$ws.Range("C7").Value = "J"

$ws.Range("A10").Value = "UNC_Fluid_Temp"
$ws.Range("C10").Value = "dTemp from ambient in deg C"
$ws.Range("B10").Value = 15

# Register the new defined names
$wb.Names.Add("UNC_Fluid_C", "=Sheet1!`$B`$9")
$wb.Names.Add("UNC_Fluid_Temp", "=Sheet1!`$B`$10")

# Widen column A to fit the longer labels (target stored width is 18.28515625
# chars; the COM ColumnWidth setter here snaps to a 6px/char grid with a 5px
# pad, i.e. stored = (round(input*6)+5)/6, so 17.5 lands on the closest
# reachable grid point, 18.33(3))
$ws.Columns.Item(1).ColumnWidth = 17.5

# Move the selection, matching the saved workbook state
$ws.Range("B15").Select() | Out-Null

$wb.Save()
